$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the original "Primera" record (currently in row 12) by moving it
# down to the newly added row 13.
$ws.Range("A13").Value = $ws.Range("A12").Value2
$ws.Range("B13").Value = $ws.Range("B12").Value2
$ws.Range("C13").Value = $ws.Range("C12").Value2
$ws.Range("D13").Value = $ws.Range("D12").Value2
$ws.Range("D13").NumberFormat = $ws.Range("D12").NumberFormat
$ws.Range("E13").Value = $ws.Range("E12").Value2
$ws.Range("F13").Value = $ws.Range("F12").Value2
$ws.Range("G13").Value = $ws.Range("G12").Value2
$ws.Range("H13").Value = $ws.Range("H12").Value2
$ws.Range("I13").Value = $ws.Range("I12").Value2
$ws.Range("J13").Value = $ws.Range("J12").Value2
$ws.Range("K13").Value = $ws.Range("K12").Value2
$ws.Range("L13").Value = $ws.Range("L12").Value2
$ws.Range("M13").Value = $ws.Range("M12").Value2
$ws.Range("N13").Value = $ws.Range("N12").Value2
$ws.Range("O13").Value = $ws.Range("O12").Value2
$ws.Range("P13").Value = $ws.Range("P12").Value2
$ws.Range("Q13").Value = $ws.Range("Q12").Value2
$ws.Range("R13").Value = $ws.Range("R12").Value2
$ws.Range("S13").Value = $ws.Range("S12").Value2
$ws.Range("T13").Value = $ws.Range("T12").Value2

# Update row 12 in place with the new "Tercera" record.
$ws.Range("D12").Value = 44783
$ws.Range("L12").Value = "Tercera"
$ws.Range("M12").Value = 100
$ws.Range("N12").Value = 27000
$ws.Range("O12").Value = 28000
$ws.Range("P12").Value = 27500
$ws.Range("S12").Value = 2292
